# Weekly update: a new week's price data (fecha 44488) is published for
# "Abate Fettel" and "Packham's Triumph" (Primera quality). The new rows
# are inserted at the top of the data block (right after the most recent
# existing row), pushing all older rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right before the current row 164, shifting
# every row from 164 downward down by two (164->166, ..., 282->284).
$ws.Rows.Item(164).Insert()
$ws.Rows.Item(164).Insert()

# Fill in the first new row (164): Abate Fettel, Primera.
$ws.Range("A164").Value = 11
$ws.Range("B164").Value = "Vega Monumental Concepción"
$ws.Range("C164").Value = "Bíobío"
$ws.Range("D164").Value = 44488
$ws.Range("E164").Value = 8
$ws.Range("F164").Value = "Fruta"
$ws.Range("G164").Value = 100104
$ws.Range("H164").Value = "Frutos de pepita"
$ws.Range("I164").Value = 100104005
$ws.Range("J164").Value = "Pera"
$ws.Range("K164").Value = "Abate Fettel"
$ws.Range("L164").Value = "Primera"
$ws.Range("M164").Value = 450
$ws.Range("N164").Value = 8000
$ws.Range("O164").Value = 9000
$ws.Range("P164").Value = 8444
$ws.Range("Q164").Value = "$/caja 17 kilos empedrada"
$ws.Range("R164").Value = "Región de O'Higgins"
$ws.Range("S164").Value = 497
$ws.Range("T164").Value = 17

# Fill in the second new row (165): Packham's Triumph, Primera.
$ws.Range("A165").Value = 11
$ws.Range("B165").Value = "Vega Monumental Concepción"
$ws.Range("C165").Value = "Bíobío"
$ws.Range("D165").Value = 44488
$ws.Range("E165").Value = 8
$ws.Range("F165").Value = "Fruta"
$ws.Range("G165").Value = 100104
$ws.Range("H165").Value = "Frutos de pepita"
$ws.Range("I165").Value = 100104005
$ws.Range("J165").Value = "Pera"
$ws.Range("K165").Value = "Packham's Triumph"
$ws.Range("L165").Value = "Primera"
$ws.Range("M165").Value = 430
$ws.Range("N165").Value = 10000
$ws.Range("O165").Value = 11000
$ws.Range("P165").Value = 10465
$ws.Range("Q165").Value = "$/caja 17 kilos empedrada"
$ws.Range("R165").Value = "Región de O'Higgins"
$ws.Range("S165").Value = 616
$ws.Range("T165").Value = 17
